# Deploy the implementation guide:
# retarget this ValueSet "include" workbook from the (HP ontology) hp.owl
# CodeSystem onto the group-member-status CodeSystem, and simplify the
# "Include" sheet from an explicit concept list to an "include all codes"
# declaration.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Metadata"
$ws2 = $wb.Worksheets.Item(2)   # "Include from hp.owl"

# --- Rename the second sheet to reflect the new source CodeSystem ---
$ws2.Name = "Include from group-member-sta"

# --- Metadata sheet: update the Description to match the Title text ---
$ws1.Range("B11").Value = "Group member status"

# --- Include sheet: drop the explicit AFF/UNF/UNK concept rows and switch
#     to a simple "All codes" declaration, then repoint the System URI at
#     the new CodeSystem. ---

# Remove the two rows that held the "UNF / Unaffected" and
# "UNK / Juvenile onset" concept rows (rows 3 & 4); delete bottom-up so
# row numbers stay valid. This leaves the blank separator row and the
# "System URI" row intact, just shifted up two rows.
$ws2.Rows.Item(4).Delete()
$ws2.Rows.Item(3).Delete()

# Row 1 used to be "Concept" / "Description" -> becomes just "Codes"
$ws2.Range("A1").Value = "Codes"
$ws2.Range("B1").Clear()

# Row 2 used to be "AFF" / "Affected" -> becomes just "All codes"
$ws2.Range("A2").Value = "All codes"
$ws2.Range("B2").Clear()

# Row 3 (formerly the blank separator row) stays blank, unchanged.

# Row 4 (formerly "System URI" / hp.owl URL) -> point at the new
# CodeSystem URL; the "System URI" label itself is unchanged.
$ws2.Range("B4").Value = "http://fhir.cqgc.ferlab.bio/CodeSystem/group-member-status"
